# "FIx Slicing & Import Excel"
#
# The "Lokasi" column is dropped from the Table4 listobject/table (and its
# header), and the little "contoh pengisian" (example-fill) reference row
# that used to live directly under the header row (A2:L2) is rebuilt further
# to the right of the table (L1:U2) so it no longer collides with the table's
# own columns/slicing. The placeholder strings "(CONTOH)1234567" /
# "JANGAN DIHAPUS/DIGANTI" are replaced by a plain "Contoh Pengisian" header
# with a realistic example row underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Table: remove the "Lokasi" column (shrinks Table4 from K to J) ---
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item("Lokasi").Delete()

# --- 2. Capture the "filled example" row's look (greenish fill, style s=1)
#        before the source cells are cleared, and stamp it onto the new
#        example block so we reuse the existing style instead of minting one.
$ws.Range("A2").Copy()
$ws.Range("L2:U2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Clear out the old in-table example-data cells (A2:K2) - this data
#        moves to the new block to the right of the table.
$ws.Range("A2:K2").Clear()

# --- 4. Build the new "Contoh Pengisian" example block beside the table ---
$ws.Range("L1").Value = "Contoh Pengisian"

$ws.Range("L2").Value = 1234567
$ws.Range("M2").Value = "Differential Pressure Transmitter"
$ws.Range("N2").Value = "0-400"
$ws.Range("O2").Value = "MMH2O"
$ws.Range("P2").Value = "DIFFERENTIAL PRESSURE TRANSMITTER 0-400 MMH2O"
$ws.Range("Q2").Value = "Baru/Bekas"
$ws.Range("R2").Value = "Yokogawa"
$ws.Range("S2").Value = "EJA110E"
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 1

# --- 5. Widen the new columns so the example block is readable ---
$ws.Columns.Item(11).ColumnWidth = 17.42578125
$ws.Columns.Item(12).ColumnWidth = 17.42578125
$ws.Columns.Item(13).ColumnWidth = 30.85546875
$ws.Columns.Item(14).ColumnWidth = 8.140625
$ws.Columns.Item(15).ColumnWidth = 8.140625
$ws.Columns.Item(16).ColumnWidth = 49
$ws.Columns.Item(17).ColumnWidth = 11
$ws.Columns.Item(18).ColumnWidth = 9.85546875
$ws.Columns.Item(19).ColumnWidth = 8

# --- 6. Move the selection onto the new example block ---
$null = $ws.Range("L4").Select()

Write-Output "Lokasi column removed; example block rebuilt at L1:U2"
